# Append extra identifier details to the "N° du prélèvement" column (B)
# for the rows that were updated during metadata model evaluation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = "24219576 1.1. Foie"
$ws.Range("B6").Value  = "24BB11466 07"
$ws.Range("B10").Value = "24MH9794 RF"
$ws.Range("B11").Value = "24MH9721 BN"
$ws.Range("B12").Value = "24EC09559 frottis 1"
$ws.Range("B14").Value = "24CU052383 pneu"
